# 1401CS75.xlsx marksheet update
# - recompute Right/Wrong/Not-Attempt/Max + Total row with numeric marking values
# - Total/Max display becomes "53/112" instead of "Absent"
# - collapses the sheet from 3 side-by-side question blocks (A:B, D:E, G:H)
#   down to a single block (A:B) plus a second block (D:E) that only goes to row 18
# - fills in "Student Ans" (col A / col D) for the surviving question(s) and
#   colors them green (correct) / red (incorrect) to match the "Correct Ans" column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: apply the formatting of an existing "template" cell (one that
# already carries the desired named style) onto a target cell without
# disturbing its value. Using Copy + PasteSpecial(formats) makes Excel reuse
# the existing style record instead of synthesising a brand-new one.
# ---------------------------------------------------------------------------
function Set-CellStyle {
    param($templateRef, $targetRef)
    $ws.Range($templateRef).Copy()
    $ws.Range($targetRef).PasteSpecial(-4122)
}

# Style templates already present in the sheet (their own style never changes):
#   mtitleStyle (bold title)  -> A9
#   correctStyle (green)      -> B10
#   incorrectStyle (red)      -> C10
#   normalStyle (plain)       -> A16
#   absoluteStyle (blue)      -> B16

# ---------------------------------------------------------------------------
# Score summary block (rows 10-12)
# ---------------------------------------------------------------------------

# Row 10 "No." - give A10 the bold mtitleStyle header look, update counts
Set-CellStyle "A9" "A10"
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

# Row 11 "Marking"
Set-CellStyle "A9" "A11"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 "Total"
Set-CellStyle "A9" "A12"
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "53/112"

# ---------------------------------------------------------------------------
# Remove the 2nd/3rd question blocks that no longer apply
# ---------------------------------------------------------------------------
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# Fill in the "Student Ans" values for the remaining question block (col D,
# rows 16-18) - all three match the "Correct Ans" column so they're green.
# ---------------------------------------------------------------------------
Set-CellStyle "B10" "D16"
$ws.Range("D16").Value = "Option A"

Set-CellStyle "B10" "D17"
$ws.Range("D17").Value = "Option C"

Set-CellStyle "B10" "D18"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------------
# Fill in "Student Ans" values for the first question block (col A,
# rows 16-40). Rows not listed keep their existing blank/normalStyle cell.
# Green (correctStyle) = matches Correct Ans, Red (incorrectStyle) = mismatch.
# ---------------------------------------------------------------------------
$correct = @{
    18 = "Option B"
    19 = "Option C"
    22 = "Option D"
    27 = "Option A"
    29 = "Option D"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
}
foreach ($row in $correct.Keys) {
    $ref = "A$row"
    Set-CellStyle "B10" $ref
    $ws.Range($ref).Value = $correct[$row]
}

$incorrect = @{
    24 = "Option B"
    26 = "Option D"
    36 = "Option B"
}
foreach ($row in $incorrect.Keys) {
    $ref = "A$row"
    Set-CellStyle "C10" $ref
    $ws.Range($ref).Value = $incorrect[$row]
}
